$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.648.09"
$ws.Range("E2").Value = "  +1.38%  "

$ws.Range("D3").Value = "3.359.37"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'585.49"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Value = "'178.31"
$ws.Range("E6").Value = "  +1.79%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +0.62%  "

$ws.Range("E9").Value = "  +4.14%  "

$ws.Range("E10").Value = "  +1.16%  "

$ws.Range("D11").Value = "'48.03"
$ws.Range("E11").Value = "  +5.91%  "

$ws.Range("E12").Value = "  +2.12%  "

$ws.Range("D13").Value = "'696.67"
$ws.Range("E13").Value = "  +5.64%  "

$ws.Range("D14").Value = "3.905.87"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").Value = "'8.49"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "68.611.58"
$ws.Range("E16").Value = "  +1.45%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.120"
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.366.39"
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").Value = "'17.54"
$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("D20").Value = "'11.27"
$ws.Range("E20").Value = "  +2.78%  "

$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("E22").Value = "  +3.85%  "

$ws.Range("D23").Value = "'16.99"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "'100.11"
$ws.Range("E24").Value = "  +1.56%  "

$ws.Range("E25").Value = "  +1.74%  "

$ws.Range("E26").Value = "  +1.80%  "

$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  +3.12%  "

$ws.Range("D28").Value = "'33.11"
$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("D29").Value = "'8.58"
$ws.Range("E29").Value = "  +1.96%  "

$ws.Range("D30").Value = "'6.99"
$ws.Range("E30").Value = "  -2.94%  "

$ws.Range("D31").Value = "'11.12"

$ws.Range("D32").Value = "'549.47"
$ws.Range("E32").Value = "  -3.32%  "

$ws.Range("E33").Value = "  +0.86%  "

$ws.Range("D34").Value = "'58.42"
$ws.Range("E34").Value = "  +3.25%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "3.708.95"
$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("D37").Value = "'3.41"
$ws.Range("E37").Value = "  +3.76%  "

$ws.Range("D38").Value = "'0.143"
$ws.Range("E38").Value = "  +9.83%  "

$ws.Range("D39").Value = "'34.70"
$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("D40").Value = "'3.20"
$ws.Range("E40").Value = "  +2.67%  "

$ws.Range("D41").Value = "'2.63"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").Value = "0.0₃0675"
$ws.Range("E42").Value = "  +2.11%  "

$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("D44").Value = "'0.0415"
$ws.Range("E44").Value = "  +2.31%  "

$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("E47").Value = "  +0.83%  "

$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").Value = "'131.90"
$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("E51").Value = "  -1.48%  "
